# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q2" and "总计",
#    carrying the per-fund holdings table for the new quarter.
# 2. Insert a new summary row at the top of "总计" (the "总计"/rollup
#    sheet) for "2022-Q1", shifting the older quarters down and
#    renumbering the running index column (A).

$wb = $excel.ActiveWorkbook

$wsQ2 = $wb.Worksheets.Item("2021-Q2")
$wsTotal = $wb.Worksheets.Item("总计")

# --- 1. New "2022-Q1" sheet, inserted right before "总计" -----------------

$newWs = $wb.Worksheets.Add($wsTotal)
$newWs.Name = "2022-Q1"

# Match the look of the other quarterly sheets: copy the header-row style
# (bold/centered/bordered) and the style used for the running-index column.
$wsQ2.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)

$wsQ2.Range("A2").Copy()
$newWs.Range("A2").PasteSpecial(-4122)

$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

$newWs.Range("A2").Value = 0

# These columns hold numeric-looking values that must stay text (as in the
# other quarterly sheets), so force text format before writing, then strip
# the leftover direct formatting so the cells match the plain data cells
# used elsewhere.
$newWs.Range("B2:G2").NumberFormat = "@"
$newWs.Range("B2").Value = "010404"
$newWs.Range("C2").Value = "博道盛利6个月持有期混合"
$newWs.Range("D2").Value = "1.29"
$newWs.Range("E2").Value = "34.13"
$newWs.Range("F2").Value = "0.50"
$newWs.Range("G2").Value = "0.0064"
$newWs.Range("B2:G2").ClearFormats()
$newWs.Range("H2").Value = 8

# --- 2. Prepend a "2022-Q1" row to "总计" ----------------------------------

# The worksheet collection shifted when the sheet above was inserted, which
# leaves the old $wsTotal handle pointing at the wrong (stale) index, so
# re-resolve it by name before touching it.
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()

# Re-apply the styles that belonged to the (now shifted) former row 2 so the
# newly-inserted row matches the sheet's look.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("B3:D3").Copy()
$wsTotal.Range("B2:D2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.01

# Renumber the running index column for the rows that shifted down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
